$wb = $excel.ActiveWorkbook

# --- Rename the first worksheet ---
# "ALZ Policy Assignments" -> "ALZ Policy Assignments 03CY23"
# (the _FilterDatabase defined name that points at this sheet is
# reference-tracked and updates automatically when the sheet is renamed)
$ws = $wb.Worksheets.Item("ALZ Policy Assignments")
$ws.Name = "ALZ Policy Assignments 03CY23"

# --- Update the view state on that sheet ---
# Scroll so row 8 is the top-left visible row, and move the active
# selection from J15 to A16.
$ws.Activate()
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
